# Renamed few transcripts. Updated the DataSheet
# Replace full speaker names in column D with their first-letter abbreviation,
# for the specific rows that were retagged (leaves other speakers, such as
# BRYAN, DANIEL, ANDREW, untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "HILLARY LEWIS-WOLFSEN" = "T"
    "STUDENT A"             = "S"
    "STUDENT B"             = "S"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 4).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 35 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()

    if ($null -ne $val -and $map.ContainsKey([string]$val)) {
        $cell.Value = $map[[string]$val]
    }
}
